$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column B (ASIN moves from B to C, etc.)
$ws.Columns.Item(2).Insert()

# New header for the inserted column
$ws.Range("B1").Value = "Week_Start_Date"

# Make sure the Week_Start_Date column is stored as text (not auto-converted to a date serial)
$ws.Range("B2:B17").NumberFormat = "@"

$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

for ($i = 0; $i -lt $weekStartDates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $weekStartDates[$i]
}

# Update the Week labels in column A to drop the leading zero (W01 -> W1 ... W09 -> W9)
$weekLabels = @(
    "W1", "W2", "W3", "W4", "W5", "W6", "W7", "W8", "W9",
    "W10", "W11", "W12", "W13", "W14", "W15", "W16"
)

for ($i = 0; $i -lt $weekLabels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weekLabels[$i]
}
